$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextCell($addr, $val) {
    $r = $ws.Range($addr)
    $r.NumberFormat = "@"
    $r.Value = $val
    $r.ClearFormats()
}

Set-TextCell "D2" "328.74"
Set-TextCell "E2" "-0.10%"
Set-TextCell "D3" "44.28"
Set-TextCell "E3" "-0.35%"
Set-TextCell "D4" "5.512"
Set-TextCell "E4" "-1.42%"
Set-TextCell "D5" "0.08104"
Set-TextCell "E5" "0.12%"
Set-TextCell "E6" "1.72%"
Set-TextCell "B7" "MXToken"
Set-TextCell "C7" "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
Set-TextCell "D7" "0.9719"
Set-TextCell "E7" "1.83%"
Set-TextCell "B8" "LiechtensteinCryptoassetsExchange"
Set-TextCell "C8" "https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx"
Set-TextCell "D8" "0.1115"
Set-TextCell "E8" "-6.35%"
Set-TextCell "B9" "WazirX"
Set-TextCell "C9" "https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx"
Set-TextCell "D9" "0.1884"
Set-TextCell "E9" "1.54%"
Set-TextCell "B10" "MCDex"
Set-TextCell "C10" "https://coinranking.com/coin/3nMM61qeg+mcdex-mcb"
Set-TextCell "D10" "10.17"
Set-TextCell "E10" "-0.43%"
Set-TextCell "B11" "MandalaExchangeToken"
Set-TextCell "C11" "https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx"
Set-TextCell "D11" "0.09962"
Set-TextCell "E11" "2.40%"
Set-TextCell "B12" "BitrueCoin"
Set-TextCell "C12" "https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr"
Set-TextCell "D12" "0.04711"
Set-TextCell "E12" "2.53%"
Set-TextCell "B13" "BitMartToken"
Set-TextCell "C13" "https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx"
Set-TextCell "D13" "0.1055"
Set-TextCell "E13" "-1.29%"
Set-TextCell "B14" "BitForexToken"
Set-TextCell "C14" "https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf"
Set-TextCell "D14" "0.001260"
Set-TextCell "E14" "-1.13%"
Set-TextCell "B15" "CoinExToken"
Set-TextCell "C15" "https://coinranking.com/coin/APDVU0XEViZ2o+coinextoken-cet"
Set-TextCell "D15" "0.04103"
Set-TextCell "E15" "-2.70%"
Set-TextCell "B16" "TigerCash"
Set-TextCell "C16" "https://coinranking.com/coin/6hIn06L2+tigercash-tch"
Set-TextCell "D16" "0.006026"
Set-TextCell "E16" "2.97%"
Set-TextCell "B17" "LEO"
Set-TextCell "C17" "https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo"
Set-TextCell "D17" "3.340"
Set-TextCell "E17" "-0.87%"
Set-TextCell "B18" "GateToken"
Set-TextCell "C18" "https://coinranking.com/coin/t7m8DZVyMsAu+gatetoken-gt"
Set-TextCell "D18" "4.424"
Set-TextCell "E18" "2.72%"
Set-TextCell "D19" "2.654"
Set-TextCell "E19" "3.48%"
Set-TextCell "D20" "0.3302"
Set-TextCell "E20" "-4.95%"
Set-TextCell "D21" "0.1389"
Set-TextCell "E21" "-1.29%"
Set-TextCell "E22" "2.70%"
Set-TextCell "D23" "0.001306"
Set-TextCell "E23" "4.69%"
Set-TextCell "D24" "0.004389"
Set-TextCell "E24" "1.50%"
Set-TextCell "E25" "7.57%"
Set-TextCell "D26" "0.0003732"
Set-TextCell "E26" "-6.20%"
Set-TextCell "D38" "0.02670"
Set-TextCell "E38" "-0.21%"
Set-TextCell "D39" "0.05642"
Set-TextCell "E39" "1.41%"
Set-TextCell "D40" "0.007605"
Set-TextCell "E40" "0.18%"
Set-TextCell "D41" "0.1412"
Set-TextCell "E41" "0.06%"
Set-TextCell "D42" "0.008229"
Set-TextCell "E42" "2.02%"
Set-TextCell "E43" "-2.84%"
Set-TextCell "D44" "0.008296"
Set-TextCell "E44" "-1.27%"
Set-TextCell "D45" "0.00007087"
Set-TextCell "E45" "-1.65%"
Set-TextCell "D46" "0.00000000748"
Set-TextCell "E46" "-0.08%"
Set-TextCell "D47" "0.0005788"
Set-TextCell "E47" "-0.41%"
Set-TextCell "D48" "0.002515"
Set-TextCell "E48" "10.93%"
Set-TextCell "D49" "0.003632"
Set-TextCell "E49" "-13.19%"
Set-TextCell "D50" "0.00002096"
Set-TextCell "E50" "-0.08%"
Set-TextCell "D51" "0.0001996"
Set-TextCell "E51" "-0.08%"
